# Adds a new "2021/12/20" row at the top of the data (row 2) on every sheet,
# pushing all the existing rows down by one. Matches the commit:
# "2021/12/20 - everyday data updated"

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 台指期換倉成本計算 (A1:F5 -> A1:F6) ---
$ws1 = $wb.Worksheets.Item("台指期換倉成本計算")
$ws1.Rows("2:2").Insert()
$ws1.Range("A2").Value = "日期：2021/12/20"
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "202202"
$ws1.Range("B2").ClearFormats()
$ws1.Range("C2").Value = 17549
$ws1.Range("D2").Value = 2109
$ws1.Range("E2").Value = 11845575
$ws1.Range("F2").Value = 17614

# --- Sheet 2: 散戶多空力道 (A1:B39 -> A1:B40) ---
$ws2 = $wb.Worksheets.Item("散戶多空力道")
$ws2.Rows("2:2").Insert()
$ws2.Range("A2").Value = "日期：2021/12/20"
$ws2.Range("B2").Value = 0.2

# --- Sheet 3: 三大法人買賣金額 (A1:C39 -> A1:C40) ---
$ws3 = $wb.Worksheets.Item("三大法人買賣金額")
$ws3.Rows("2:2").Insert()
$ws3.Range("A2").Value = "110年12月20日"
$ws3.Range("B2").Value = -194.23
$ws3.Range("C2").Value = 50.75

# --- Sheet 4: 大盤多空點位 (A1:B38 -> A1:B39) ---
$ws4 = $wb.Worksheets.Item("大盤多空點位")
$ws4.Rows("2:2").Insert()
$ws4.Range("A2").Value = "110年12月20日"
$ws4.Range("B2").Value = 17720.87

# --- Sheet 5: 期貨大額交易人未沖銷部位 (A1:N37 -> A1:N38) ---
$ws5 = $wb.Worksheets.Item("期貨大額交易人未沖銷部位")
$ws5.Rows("2:2").Insert()
$ws5.Range("A2").NumberFormat = "@"
$ws5.Range("A2").Value = "2021/12/20"
$ws5.Range("A2").ClearFormats()
$ws5.Range("B2").Value = 46944
$ws5.Range("C2").Value = 55396
$ws5.Range("D2").Value = 1101
$ws5.Range("E2").Value = 1087
$ws5.Range("F2").Value = 21589
$ws5.Range("G2").Value = 48231
$ws5.Range("H2").Value = -426
$ws5.Range("I2").Value = 802
$ws5.Range("J2").Value = -26642
$ws5.Range("K2").Value = -1228
$ws5.Range("L2").Value = 1527
$ws5.Range("M2").Value = 285
$ws5.Range("N2").Value = 1242

Write-Output "Inserted 2021/12/20 rows on all 5 sheets"
